# [FE-Edit] Download Product Template from Backend
#
# Add a new "existing_strategy" reference sheet (after "existing_coa") that
# lists the available strategy names, mirroring the layout/formatting of the
# other "existing_*" lookup sheets already in the template.

$wb = $excel.ActiveWorkbook

$coaSheet = $wb.Worksheets.Item("existing_coa")

# Insert the new sheet immediately after "existing_coa" so the tab order
# becomes: product, example, existing_coa, existing_strategy.
$newSheet = $wb.Worksheets.Add($null, $coaSheet)
$newSheet.Name = "existing_strategy"

# Header row
$newSheet.Range("A1").Value = "strategy_name"

# Data rows
$newSheet.Range("A2").Value = "Grow the Business"
$newSheet.Range("A3").Value = "Strategy 1"
$newSheet.Range("A4").Value = "Run the Business"
$newSheet.Range("A5").Value = "Strategy 2"

# Match the header styling used by the other lookup sheets (e.g.
# "existing_coa" A1): bold font, thin border all around, centered alignment.
# Copy the format straight from that header cell (format-painter style) so
# the existing style is reused rather than a brand new one being created.
$coaSheet.Range("A1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
